# The document has a single paragraph containing "каккакакак" followed
# by the (hidden) "_GoBack" bookmark. The edit:
#   1. appends a new run "zxzxzxzxzxzxxzx" right before the bookmark,
#   2. tags that new run as English (US): <w:rPr><w:lang w:val="en-US"/></w:rPr>
#   3. tags the paragraph mark itself as English (US) too (Word's "word
#      correctives" updated the paragraph mark language along with the
#      freshly typed run): <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>

$d = $word.ActiveDocument

$newText = "zxzxzxzxzxzxxzx"
$newLang = "en-US"

$para = $d.Paragraphs(1)

# Grab the paragraph's existing text content, excluding the trailing
# paragraph mark, so the original run is carried over unchanged.
$paraRange = $para.Range
$existingText = $d.Range($paraRange.Start, $paraRange.End - 1).Text

function Escape-Xml([string]$s) {
    $s = $s.Replace("&", "&amp;")
    $s = $s.Replace("<", "&lt;")
    $s = $s.Replace(">", "&gt;")
    $s = $s.Replace('"', "&quot;")
    return $s
}

$existingTextXml = Escape-Xml($existingText)
$newTextXml = Escape-Xml($newText)

# Rebuild the whole (single) paragraph's content in one shot so the new
# run lands before the bookmark and the paragraph-mark language sticks,
# matching what real Word does when the user types at the end of the line.
$bodyXml = "<w:p>" + `
    "<w:pPr><w:rPr><w:lang w:val=`"$newLang`"/></w:rPr></w:pPr>" + `
    "<w:r><w:t>$existingTextXml</w:t></w:r>" + `
    "<w:r><w:rPr><w:lang w:val=`"$newLang`"/></w:rPr><w:t>$newTextXml</w:t></w:r>" + `
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" + `
    "<w:bookmarkEnd w:id=`"0`"/>" + `
    "</w:p>"

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    "<w:body>$bodyXml</w:body>" + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($packageXml)

$d.Save()
